$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "statement" column (L) duplicated the <statement> row header and was
# removed from the LL parsing table; deleting the whole column shifts the
# remaining columns (M..P) left by one and drops the now-unused
# "statement" shared string automatically.
$ws.Range("L:L").Delete()

# Fill in the newly documented table entries in row 2 (<program> row).
$ws.Range("D2").Value = 25
$ws.Range("H2").Value = 24
$ws.Range("I2").Value = 26
$ws.Range("J2").Value = 27
$ws.Range("K2").Value = 28

# Row 11 (<idwhat>) gains a third entry.
$ws.Range("N11").Value = 32

# Row 12 (<assign>) values were bumped by one.
$ws.Range("D12").Value = 33
$ws.Range("H12").Value = 34

# Restore the view state (zoom + active cell) recorded for the sheet.
[void]$ws.Range("J8").Select()
$excel.ActiveWindow.Zoom = 130
